# Chore(Routes): Updating PUT - POST Routes Configuration
#
# The "Routes Flashcards" table used a single placeholder URL "/cards/:id"
# for the "Voir une flashcard" (GET), "Modifier une flashcard" (PUT) and
# "Supprimer une flashcard" (DELETE) routes. Update them to the real,
# deck-scoped route "/decks/:deckId/cards/:cardId" to match how the API
# is actually mounted (mirrors the existing "/decks/:deckId/cards" routes
# right above them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C29").Value = "/decks/:deckId/cards/:cardId"
$ws.Range("C31").Value = "/decks/:deckId/cards/:cardId"
$ws.Range("C32").Value = "/decks/:deckId/cards/:cardId"

# Move the active selection to where the author left off editing.
$ws.Range("E23").Select()
